# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta" edit:
#  - Moves the "JOSE ALFREDO LOZANO ZABALETA" record (previously row 18) to the top of the
#    worker list (row 16), shifting the two preceding records (AURA ESTHER FLOREZ QUIROZ,
#    CARLOS ARTURO RAMOS SUAREZ) down by one row each.
#  - Updates two "Salario Basico" values further down the table:
#      EDGAR ADOLFO CASTRO GENES  : 0       -> 828116
#      LUIS ENRIQUE CARABALLO LOPEZ: 2352990 -> 2219802

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data columns used by the table: B (Tipo Doc), C (N Doc), D (Nombre), E (Periodo Mora),
# F (Valor Mora), G (Salario Basico). Worker rows run from 16 to 29.
$firstRow = 16
$lastRow = 29
$cols = @("B", "C", "D", "E", "F", "G")

# 1) Snapshot all current worker rows (values only) before mutating anything.
$rows = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $rows[$r] = $rowVals
}

# 2) Build the new row order: the record that used to sit on row 18 (JOSE ALFREDO LOZANO
#    ZABALETA) moves to row 16; the records that used to occupy rows 16 and 17 shift down
#    to rows 17 and 18. Rows 19-29 stay in the same relative order.
$newOrder = @(18, 16, 17)
for ($r = 19; $r -le $lastRow; $r++) {
    $newOrder += $r
}

# 3) Write the reordered data back out, row by row.
$destRow = $firstRow
foreach ($srcRow in $newOrder) {
    $srcVals = $rows[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
    $destRow++
}

# 4) Apply the two updated "Salario Basico" (column G) values, matched by worker name so the
#    edit is robust regardless of which physical row each worker now occupies.
$updates = @{
    "EDGAR ADOLFO CASTRO GENES"    = 828116
    "LUIS ENRIQUE CARABALLO LOPEZ" = 2219802
}

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $name = $ws.Range("D$r").Value2
    if ($updates.ContainsKey($name)) {
        $ws.Range("G$r").Value = $updates[$name]
    }
}
